# Generate Report for Handback
# Update the "Latest Handoff Datetime" (D) and "Latest Handback DateTime" (G)
# values for the e63996ec-... file row (row 3) on the per-language handback
# report sheets, reflecting the newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-03-11 01:08:07"
$wsZhCn.Range("G3").Value = "2016-03-11 01:08:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-03-11 01:08:14"
$wsDeDe.Range("G3").Value = "2016-03-11 01:08:57"
